# Daily attendance processing - reorders the "Recorded By" (column G) list
# so that entries get rebuilt with the literal token "System" reflected at
# the front of the list when present, by reversing the comma-separated
# list of recorders. Lists that do not contain the exact token "System"
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ','
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts -contains 'System') {
        $reversed = $parts[($parts.Length - 1)..0]
        $newVal = [string]::Join(', ', $reversed)
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
